$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('展览')
$ws.Range('F2').Value = 1327
$ws.Range('F6').Value = 37
$ws.Range('C7').Value = '上海·魔都劳动节漫展-CF01'
$ws.Range('D7').Value = '澳门路168号 月星家居（澳门路）'
$ws.Range('E7').Value = '2024.05.01 10:00-05.05 16:00'
$ws.Range('F7').Value = 778
$ws.Range('G7').Value = 59
$ws.Range('H7').Value = 'https://show.bilibili.com/platform/detail.html?id=82992'
$ws.Range('I7').Value = '//i2.hdslb.com/bfs/openplatform/202403/I7O9LMtb1710752670542.jpeg'
$ws.Range('B8').Value = '2024-05-02'
$ws.Range('C8').Value = '上海·2024GAF插画艺术节'
$ws.Range('D8').Value = '博成路850号 上海世博展览馆'
$ws.Range('E8').Value = '2024.05.02 10:30-05.04 19:00'
$ws.Range('F8').Value = 1483
$ws.Range('G8').Value = 128
$ws.Range('H8').Value = 'https://show.bilibili.com/platform/detail.html?id=83699'
$ws.Range('I8').Value = '//i1.hdslb.com/bfs/openplatform/202403/APlNld8y1711825700811.jpeg'
$ws.Range('C9').Value = '上海·「星铁LAND」2024星穹铁道嘉年华'
$ws.Range('D9').Value = '崧泽大道333号 上海国家会展中心'
$ws.Range('E9').Value = '2024.05.02 09:00-05.04 17:00'
$ws.Range('F9').Value = 97353
$ws.Range('G9').Value = '不可售'
$ws.Range('H9').Value = 'https://show.bilibili.com/platform/detail.html?id=84096'
$ws.Range('I9').Value = '//i1.hdslb.com/bfs/openplatform/202404/CW93VZON1712826642232.jpeg'
$ws.Range('C10').Value = '上海·第五十八届燃梦星辰国潮嘉年华-随机宅舞'
$ws.Range('D10').Value = '周家嘴路3608号 宝龙旭辉广场'
$ws.Range('E10').Value = '2024.05.02 10:20-05.03 16:30'
$ws.Range('F10').Value = 806
$ws.Range('G10').Value = 68
$ws.Range('H10').Value = 'https://show.bilibili.com/platform/detail.html?id=82761'
$ws.Range('I10').Value = '//i0.hdslb.com/bfs/openplatform/202403/azEA4EM01710236719279.jpeg'
$ws.Range('B11').Value = '2024-05-03'
$ws.Range('C11').Value = '上海·2024明日方舟嘉年华'
$ws.Range('D11').Value = '崧泽大道333号 上海国家会展中心'
$ws.Range('E11').Value = '2024.05.03 09:00-05.05 18:00'
$ws.Range('F11').Value = 40484
$ws.Range('G11').Value = '不可售'
$ws.Range('H11').Value = 'https://show.bilibili.com/platform/detail.html?id=83707'
$ws.Range('I11').Value = '//i1.hdslb.com/bfs/openplatform/202404/QkxIHGSy1712110232653.jpeg'
$ws.Range('F12').Value = 814
$ws.Range('F13').Value = 96
$ws.Range('F14').Value = 610
$ws.Range('F16').Value = 710
$ws.Range('F17').Value = 1359
$ws.Range('F18').Value = 226
$ws.Range('F20').Value = 184
$ws.Range('F22').Value = 5398
$ws.Range('F23').Value = 312
$ws.Range('F25').Value = 2507
$ws.Range('F26').Value = 5989
$ws.Range('F27').Value = 139
$ws.Range('F28').Value = 1050
$ws.Range('F29').Value = 623
$ws.Range('F30').Value = 70
$ws.Range('F32').Value = 1069
$ws.Range('F33').Value = 34
$ws.Range('F35').Value = 87
$ws.Range('F37').Value = 745
$ws.Range('F39').Value = 57
$ws.Range('F41').Value = 1096
$ws.Range('F44').Value = 69
$ws.Range('F45').Value = 29
$ws.Range('F46').Value = 122
$ws.Range('F47').Value = 644
$ws.Range('F48').Value = 16
$ws.Range('F49').Value = 35
$ws = $wb.Worksheets.Item('演出')
$ws.Range('F5').Value = 2099
$ws.Range('G5').Value = '不可售'
$ws.Range('F11').Value = 693
$ws.Range('F12').Value = 21
$ws.Range('F24').Value = 533
$ws.Range('F29').Value = 109
$ws.Range('F32').Value = 73
$ws.Range('F36').Value = 225
$ws.Range('F37').Value = 912
$ws.Range('F38').Value = 511
$ws.Range('F40').Value = 43
$ws.Range('F43').Value = 81
$ws.Range('F46').Value = 22
$ws = $wb.Worksheets.Item('本地生活')
$ws.Range('F5').Value = 792
$ws.Range('F6').Value = 456
$ws.Range('F7').Value = 255
$ws = $wb.Worksheets.Item('全部类型')
$ws.Range('F3').Value = 1327
$ws.Range('F6').Value = 456
$ws.Range('F7').Value = 255
$ws.Range('F8').Value = 255
$ws.Range('F10').Value = 37
$ws.Range('F12').Value = 1483
$ws.Range('F13').Value = 806
$ws.Range('F14').Value = 96
$ws.Range('F15').Value = 693
$ws.Range('F16').Value = 1359
$ws.Range('F17').Value = 226
$ws.Range('F20').Value = 184
$ws.Range('F21').Value = 312
$ws.Range('F22').Value = 2507
$ws.Range('F23').Value = 5989
$ws.Range('F24').Value = 139
$ws.Range('F25').Value = 1050
$ws.Range('F28').Value = 623
$ws.Range('F29').Value = 70
$ws.Range('F30').Value = 1069
$ws.Range('F32').Value = 87
$ws.Range('F34').Value = 745
$ws.Range('F35').Value = 73
$ws.Range('F36').Value = 57
$ws.Range('F37').Value = 1096
$ws.Range('F43').Value = 43
$ws.Range('F44').Value = 122
$ws.Range('F46').Value = 81
$ws.Range('F49').Value = 35
$ws.Range('F50').Value = 22
